$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    if ($row -ge 50 -and $row -le 57) {
        $f = $ws.Cells.Item($row, 6).Value2
        $ws.Cells.Item($row, 5).Value = 10
        $ws.Cells.Item($row, 6).Value = $f + 10
    } else {
        $e = $ws.Cells.Item($row, 5).Value2
        $ws.Cells.Item($row, 5).Value = $e - 1
    }
}
